# "updating TMR and DR" - Requirements Traceability Matrix update.
#
# Navigation test rows (B9:B14, "Navigation" requirement) get their
# "All links viable #1" placeholder text disambiguated into #2/#3/#4 for
# rows 10-12, and the two remaining "Back Navigation" / "Home Page Links
# Visible" rows (13/14) get their Test Details (E), defects note (G) and
# defectID count (H) filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Disambiguate the duplicated "All links viable #1" test descriptions ---
$ws.Range("D10").Value = "All links viable #2"
$ws.Range("D11").Value = "All links viable #3"
$ws.Range("D12").Value = "All links viable #4"

# --- Row 13 (Back Navigation): fill in the Test Details ---
$ws.Range("E13").Value = "The manager clicks on Matrices, the title of the page should be matrix page, the manager clicks back browser button and should be on home page with the title home, The manager clicks on test cases then the browser back button again to go back to home page where the title should be home."

# --- Row 14 defect note, Row 14 Test Details, Row 13 defect note ---
$ws.Range("G14").Value = "Expected Defect Reporting but actual is Report a Defect"
$ws.Range("E14").Value = "The manager should see links for Matrices, Test Cases, Defect Reporting and Defect Overview"
$ws.Range("G13").Value = "Expected Matrix Page but actual is Matrix Dashboard"

# --- New defectID counts for the two failed Navigation test cases ---
$ws.Range("H13").Value = 2
$ws.Range("H14").Value = 1

# --- Widen the Test Description / Test Details columns to fit the new text ---
$ws.Columns.Item(4).ColumnWidth = 41.83
$ws.Columns.Item(5).ColumnWidth = 33.17

# --- Selection moves to the last-edited cell ---
$ws.Range("G14").Select()
